# Remove the "This is a sentence 5 modified." paragraph and the
# trailing empty paragraph that followed it (testing content cleanup).
$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$target = $null

for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "This is a sentence 5*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $startDelete = $target.Range.Start
    # Extend the deletion through the very end of the document content so
    # any empty paragraph(s) left trailing after the target are removed too.
    $endDelete = $d.Content.End

    $r = $d.Range($startDelete, $endDelete)
    $r.Delete()
}
